# Applies corrected IFRS financial figures for rows 2-9 (company_list sheet).
# The prior figures were off by roughly two orders of magnitude (unit/scale error);
# this restates columns D:AJ per row with the corrected values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$rowUpdates = [ordered]@{
    "D2" = 13594
    "E2" = 597
    "F2" = 597
    "G2" = 1265
    "H2" = 1012
    "I2" = 1006
    "J2" = 6
    "K2" = 13794
    "L2" = 6013
    "M2" = 7781
    "N2" = 7691
    "O2" = 90
    "P2" = 169
    "Q2" = 907
    "R2" = -789
    "S2" = -165
    "T2" = 635
    "U2" = 272
    "V2" = 2044
    "W2" = 4.39
    "X2" = 7.45
    "Y2" = 13.87
    "Z2" = 7.79
    "AA2" = 77.27
    "AB2" = 4442.64
    "AC2" = 2971
    "AD2" = 5.86
    "AE2" = 22711
    "AF2" = 0.77
    "AG2" = 200
    "AH2" = 1.15
    "AI2" = 6.73
    "AJ2" = 33865090
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 3
$rowUpdates = [ordered]@{
    "D3" = 13951
    "E3" = 633
    "F3" = 633
    "G3" = 1102
    "H3" = 727
    "I3" = 705
    "J3" = 22
    "K3" = 15483
    "L3" = 7027
    "M3" = 8456
    "N3" = 8342
    "O3" = 113
    "P3" = 169
    "Q3" = 1542
    "R3" = -1488
    "S3" = 506
    "T3" = 1276
    "U3" = 266
    "V3" = 2735
    "W3" = 4.54
    "X3" = 5.21
    "Y3" = 8.789999999999999
    "Z3" = 4.97
    "AA3" = 83.09999999999999
    "AB3" = 4818.54
    "AC3" = 2082
    "AD3" = 8.289999999999999
    "AE3" = 24634
    "AF3" = 0.7
    "AG3" = 260
    "AH3" = 1.51
    "AI3" = 12.49
    "AJ3" = 33865090
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 4
$rowUpdates = [ordered]@{
    "D4" = 16192
    "E4" = 993
    "F4" = 993
    "G4" = 1538
    "H4" = 1148
    "I4" = 1137
    "J4" = 11
    "K4" = 16323
    "L4" = 6832
    "M4" = 9492
    "N4" = 9367
    "O4" = 125
    "P4" = 169
    "Q4" = 972
    "R4" = -894
    "S4" = -206
    "T4" = 716
    "U4" = 256
    "V4" = 2722
    "W4" = 6.13
    "X4" = 7.09
    "Y4" = 12.85
    "Z4" = 7.22
    "AA4" = 71.97
    "AB4" = 5416.62
    "AC4" = 3359
    "AD4" = 6.42
    "AE4" = 27661
    "AF4" = 0.78
    "AG4" = 260
    "AH4" = 1.21
    "AI4" = 7.74
    "AJ4" = 33865090
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 5
$rowUpdates = [ordered]@{
    "D5" = 14877
    "E5" = 635
    "F5" = 635
    "G5" = 1181
    "H5" = 979
    "I5" = 956
    "J5" = 11
    "K5" = 15914
    "L5" = 5761
    "M5" = 10153
    "N5" = 10020
    "O5" = 121
    "P5" = 169
    "Q5" = 1512
    "R5" = -1238
    "S5" = -353
    "T5" = 641
    "U5" = 872
    "V5" = 2202
    "W5" = 4.27
    "X5" = 6.58
    "Y5" = 9.869999999999999
    "Z5" = 6.07
    "AA5" = 56.75
    "AB5" = 5923.63
    "AC5" = 2824
    "AD5" = 8.68
    "AE5" = 29589
    "AF5" = 0.83
    "AG5" = 400
    "AH5" = 1.63
    "AI5" = 14.17
    "AJ5" = 33865090
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 6
$rowUpdates = [ordered]@{
    "D6" = 16020
    "E6" = 167
    "F6" = 167
    "G6" = 421
    "H6" = 318
    "I6" = 301
    "K6" = 17095
    "L6" = 6877
    "M6" = 10218
    "N6" = 10048
    "P6" = 169
    "Q6" = 571
    "R6" = -630
    "S6" = -38
    "T6" = 756
    "U6" = -185
    "V6" = 2346
    "W6" = 1.04
    "X6" = 1.98
    "Y6" = 3
    "Z6" = 1.93
    "AA6" = 67.3
    "AB6" = 5966.26
    "AC6" = 889
    "AD6" = 22.45
    "AE6" = 29671
    "AF6" = 0.67
    "AG6" = 400
    "AH6" = 2.01
    "AI6" = 45.01
    "AJ6" = 33865090
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 7
$rowUpdates = [ordered]@{
    "D7" = 22714
    "E7" = 633
    "G7" = 1173
    "H7" = 957
    "I7" = 925
    "K7" = 21035
    "L7" = 8397
    "M7" = 12638
    "N7" = 12446
    "P7" = 226
    "Q7" = -60
    "R7" = -644
    "S7" = 1887
    "T7" = 1112
    "U7" = -550
    "W7" = 2.79
    "X7" = 4.21
    "Y7" = 8.23
    "Z7" = 5.02
    "AA7" = 66.45
    "AC7" = 2071
    "AD7" = 8.09
    "AE7" = 26796
    "AF7" = 0.63
    "AG7" = 400
    "AH7" = 2.39
    "AI7" = 20.84
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 8
$rowUpdates = [ordered]@{
    "D8" = 27094
    "E8" = 1121
    "G8" = 1603
    "H8" = 1214
    "I8" = 1175
    "K8" = 22645
    "L8" = 9009
    "M8" = 13636
    "N8" = 13412
    "P8" = 226
    "Q8" = 1872
    "R8" = -1203
    "S8" = -185
    "T8" = 788
    "U8" = 818
    "W8" = 4.14
    "X8" = 4.48
    "Y8" = 9.09
    "Z8" = 5.56
    "AA8" = 66.06
    "AC8" = 2437
    "AD8" = 6.87
    "AE8" = 28874
    "AF8" = 0.58
    "AG8" = 416
    "AH8" = 2.48
    "AI8" = 17.07
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

# Row 9
$rowUpdates = [ordered]@{
    "D9" = 29006
    "E9" = 1336
    "G9" = 1777
    "H9" = 1348
    "I9" = 1303
    "K9" = 24181
    "L9" = 9419
    "M9" = 14766
    "N9" = 14502
    "P9" = 226
    "Q9" = 1667
    "R9" = -1187
    "S9" = -197
    "T9" = 806
    "U9" = 931
    "W9" = 4.61
    "X9" = 4.65
    "Y9" = 9.34
    "Z9" = 5.76
    "AA9" = 63.79
    "AC9" = 2704
    "AD9" = 6.2
    "AE9" = 31222
    "AF9" = 0.54
    "AG9" = 420
    "AH9" = 2.51
    "AI9" = 15.53
}
foreach ($addr in $rowUpdates.Keys) {
    $ws.Range($addr).Value = $rowUpdates[$addr]
}

Write-Output "Updated 244 cells across rows 2-9."
